$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1407.6666
$ws.Range("I2").Value = 1493.6666
$ws.Range("K2").Value = 1493.6666
$ws.Range("M2").Value = -1380.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1362.3235
$ws.Range("I15").Value = 1362.3235
$ws.Range("K15").Value = 4086.9705
$ws.Range("M15").Value = -3917.9705

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6950
$ws.Range("J64").Value = 9000
$ws.Range("L64").Value = 9000
$ws.Range("N64").Value = -9496

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 6950
$ws.Range("J67").Value = 9000
$ws.Range("L67").Value = 9000
$ws.Range("N67").Value = -10716

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 1428.091
$ws.Range("I82").Value = 1428.091
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 4284.272999999999
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -3878.272999999999
$ws.Range("N82").Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 1428.091
$ws.Range("I85").Value = 1428.091
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 4284.272999999999
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2880.272999999999
$ws.Range("N85").Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 2999.6667
$ws.Range("J97").Value = 2999.6667
$ws.Range("L97").Value = 8999.000100000001
$ws.Range("N97").Value = -9991.000100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 41668652
$ws.Range("J101").Value = 5192.5
$ws.Range("L101").Value = 15577.5
$ws.Range("N101").Value = -18821.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1130.4445
$ws.Range("I103").Value = 935
$ws.Range("K103").Value = 2805
$ws.Range("M103").Value = -2219

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 333333340
$ws.Range("I107").Value = 333333340
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 333333340
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -333331420
$ws.Range("N107").Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 8550775
$ws.Range("J111").Value = 966.6667
$ws.Range("L111").Value = 2900.0001
$ws.Range("N111").Value = -9034.000100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 6788.7393
$ws.Range("I112").Value = 2187.5
$ws.Range("J112").Value = 7226.952
$ws.Range("K112").Value = 6562.5
$ws.Range("L112").Value = 21680.856
$ws.Range("M112").Value = -5454.5
$ws.Range("N112").Value = -23896.856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 50002692
$ws.Range("I129").Value = 100001270
$ws.Range("K129").Value = 300003810
$ws.Range("M129").Value = -299998810

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 892.0454999999999
$ws.Range("I135").Value = 862.94116
$ws.Range("K135").Value = 7766.47044
$ws.Range("M135").Value = -5231.47044

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6370.3335
$ws.Range("I63").Value = 1833.25
$ws.Range("K63").Value = 1833.25
$ws.Range("M63").Value = -1147.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 6370.3335
$ws.Range("I66").Value = 1833.25
$ws.Range("K66").Value = 9166.25
$ws.Range("M66").Value = -5734.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 130122.6
$ws.Range("I74").Value = 13389.615
$ws.Range("J74").Value = 888887
$ws.Range("K74").Value = 13389.615
$ws.Range("L74").Value = 888887
$ws.Range("M74").Value = -12515.615
$ws.Range("N74").Value = -890635

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 130122.6
$ws.Range("I77").Value = 13389.615
$ws.Range("J77").Value = 888887
$ws.Range("K77").Value = 66948.075
$ws.Range("L77").Value = 4444435
$ws.Range("M77").Value = -62580.075
$ws.Range("N77").Value = -4453171

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 10836.471
$ws.Range("I132").Value = 13373.6
$ws.Range("J132").Value = 7212
$ws.Range("K132").Value = 40120.8
$ws.Range("L132").Value = 21636
$ws.Range("M132").Value = -37590.8
$ws.Range("N132").Value = -26696

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2393198.5
$ws.Range("I94").Value = 3135429.8
$ws.Range("J94").Value = 1564.5555
$ws.Range("K94").Value = 3135429.8
$ws.Range("L94").Value = 1564.5555
$ws.Range("M94").Value = -3134978.8
$ws.Range("N94").Value = -2466.5555

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 59944.5
$ws.Range("J109").Value = 59944.5
$ws.Range("L109").Value = 59944.5
$ws.Range("N109").Value = -62718.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 26340.2
$ws.Range("J110").Value = 26340.2
$ws.Range("L110").Value = 26340.2
$ws.Range("N110").Value = -34520.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 20371.066
$ws.Range("I134").Value = 23857.2
$ws.Range("K134").Value = 71571.60000000001
$ws.Range("M134").Value = -69036.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15361.2
$ws.Range("I31").Value = 8430.134
$ws.Range("J31").Value = 16747.414
$ws.Range("K31").Value = 8430.134
$ws.Range("L31").Value = 16747.414
$ws.Range("M31").Value = -8135.134
$ws.Range("N31").Value = -17337.414

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 15361.2
$ws.Range("I34").Value = 8430.134
$ws.Range("J34").Value = 16747.414
$ws.Range("K34").Value = 8430.134
$ws.Range("L34").Value = 16747.414
$ws.Range("M34").Value = -8228.134
$ws.Range("N34").Value = -17151.414

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 30375
$ws.Range("J59").Value = 45750
$ws.Range("L59").Value = 45750
$ws.Range("N59").Value = -48040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11764.104
$ws.Range("I86").Value = 9865
$ws.Range("K86").Value = 9865
$ws.Range("M86").Value = -8742

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 11764.104
$ws.Range("I89").Value = 9865
$ws.Range("K89").Value = 49325
$ws.Range("M89").Value = -43709

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 941.25
$ws.Range("I94").Value = 750.4
$ws.Range("J94").Value = 1028
$ws.Range("K94").Value = 750.4
$ws.Range("L94").Value = 1028
$ws.Range("M94").Value = -299.4
$ws.Range("N94").Value = -1930

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 94029.82000000001
$ws.Range("I132").Value = 113925.445
$ws.Range("K132").Value = 341776.335
$ws.Range("M132").Value = -339246.335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 8070.5137
$ws.Range("I134").Value = 5767.885
$ws.Range("K134").Value = 17303.655
$ws.Range("M134").Value = -14768.655

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1642
$ws.Range("I107").Value = 975
$ws.Range("J107").Value = 2023.1428
$ws.Range("K107").Value = 2925
$ws.Range("L107").Value = 6069.428400000001
$ws.Range("M107").Value = -1005
$ws.Range("N107").Value = -9909.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1703.5385
$ws.Range("I132").Value = 1050.2858
$ws.Range("K132").Value = 9452.572200000001
$ws.Range("M132").Value = -6922.572200000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 13880
$ws.Range("I38").Value = 7500
$ws.Range("J38").Value = 15475
$ws.Range("K38").Value = 7500
$ws.Range("L38").Value = 15475
$ws.Range("M38").Value = -7037
$ws.Range("N38").Value = -16401

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 29821.5
$ws.Range("J98").Value = 29821.5
$ws.Range("L98").Value = 29821.5
$ws.Range("N98").Value = -35811.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 994072.4399999999
$ws.Range("I122").Value = 1117831.5
$ws.Range("K122").Value = 3353494.5
$ws.Range("M122").Value = -3351044.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8285.429
$ws.Range("I7").Value = 5666.6665
$ws.Range("K7").Value = 5666.6665
$ws.Range("M7").Value = -5554.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 51490.25
$ws.Range("J22").Value = 2799.6365
$ws.Range("L22").Value = 2799.6365
$ws.Range("N22").Value = -3389.6365

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 51490.25
$ws.Range("J27").Value = 2799.6365
$ws.Range("L27").Value = 2799.6365
$ws.Range("N27").Value = -3013.6365

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 24800
$ws.Range("J95").Value = 24800
$ws.Range("L95").Value = 24800
$ws.Range("N95").Value = -30292

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4332.8096
$ws.Range("J100").Value = 6331.3335
$ws.Range("L100").Value = 6331.3335
$ws.Range("N100").Value = -7413.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 18083.834
$ws.Range("J103").Value = 18083.834
$ws.Range("L103").Value = 18083.834
$ws.Range("N103").Value = -20427.834

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 6184.75
$ws.Range("J106").Value = 6184.75
$ws.Range("L106").Value = 6184.75
$ws.Range("N106").Value = -8708.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 8285.429
$ws.Range("I126").Value = 5666.6665
$ws.Range("K126").Value = 16999.9995
$ws.Range("M126").Value = -14529.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 14612.0625
$ws.Range("I132").Value = 16092.429
$ws.Range("K132").Value = 48277.287
$ws.Range("M132").Value = -45747.287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 109949
$ws.Range("J27").Value = 109949
$ws.Range("L27").Value = 109949
$ws.Range("N27").Value = -110087

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 60000
$ws.Range("J105").Value = 60000
$ws.Range("L105").Value = 60000
$ws.Range("N105").Value = -66988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 61146
$ws.Range("J123").Value = 61146
$ws.Range("L123").Value = 61146
$ws.Range("N123").Value = -70946

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3484.4285
$ws.Range("I126").Value = 4248.25
$ws.Range("J126").Value = 2466
$ws.Range("K126").Value = 12744.75
$ws.Range("L126").Value = 7398
$ws.Range("M126").Value = -10274.75
$ws.Range("N126").Value = -12338
